# Update existing requirement rows (2 and 3) with revised copy
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Track leave requests"
$ws.Range("B2").Value = "Employees don’t know request status"
$ws.Range("C2").Value = "Add real-time tracking in the HR portal"

$ws.Range("A3").Value = "Centralize approvals"
$ws.Range("B3").Value = "Managers respond via scattered emails"
$ws.Range("C3").Value = "Create a centralized approval dashboard"

# New sample/test requirement rows (4-9)
$ws.Range("A4").Value = "Automate report generation"
$ws.Range("B4").Value = "Manual monthly reporting is slow"
$ws.Range("C4").Value = "Auto-generate PDF reports from database"

$ws.Range("A5").Value = "Improve data search"
$ws.Range("B5").Value = "Filtering records is time-consuming"
$ws.Range("C5").Value = "Add a search and filter panel for users"

$ws.Range("A6").Value = "Analyze process delays"
$ws.Range("B6").Value = "No visibility into bottlenecks"
$ws.Range("C6").Value = "Include timestamp logging and delay metrics"

$ws.Range("A7").Value = "Summarize user feedback"
$ws.Range("B7").Value = "Unstructured notes are hard to review"
$ws.Range("C7").Value = "Store feedback with sentiment tagging"

$ws.Range("A8").Value = "Visualize data trends"
$ws.Range("B8").Value = "Stakeholders don’t understand raw data"
$ws.Range("C8").Value = "Use bar charts and KPIs in reports"

$ws.Range("A9").Value = "Reduce duplicate requests"
$ws.Range("B9").Value = "Employees submit same requests twice"
$ws.Range("C9").Value = "Add validation for existing entries"

# Widen the new "Proposed Solution" column to fit the longer text
$ws.Columns.Item(3).ColumnWidth = 35.5

# Re-assert the header formatting so the bold font entry gets normalized
$ws.Range("A1:C1").Font.Bold = $true

# Match the author's view state: zoomed to 115% with C1 selected
$excel.ActiveWindow.Zoom = 115
$ws.Range("C1").Select() | Out-Null
